$d = $word.ActiveDocument
$lb = [char]11   # Word's internal "line break" character, renders as <w:br/>

# ---------------------------------------------------------------------------
# Block 1: "Programa" section, English (italic) paragraph -> split into 5
# numbered sentences separated by line breaks.
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute(
    "1.Generic biotechnological process: schematic representation; description of the main stages.2.Equipment sterilization: terminology; sterilization by physical agents; sterilization by chemical agents.3.Media sterilization by steam heating: kinetics of thermal destruction of microorganisms; destruction of nutrient media; calculation of sterilization time for batch processes; design of sterilization systems for continuous processes.4.Sterilization by filtration: microbial aerosols; air samplers; filter sizing; media sterilization.5.Kinetics and stoichiometry of microbial growth and product formation: transformation rates and conversion factors; classification of fermentative processes based on cell growth and product formation rates; influence of substrate concentration on cell growth rate; stoichiometry of microbial growth and product formation.",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0
) | Out-Null
$part1 = @(
    "1.Generic biotechnological process: schematic representation; description of the main stages.",
    "2.Equipment sterilization: terminology; sterilization by physical agents; sterilization by chemical agents.",
    "3.Media sterilization by steam heating: kinetics of thermal destruction of microorganisms; destruction of nutrient media; calculation of sterilization time for batch processes; design of sterilization systems for continuous processes.",
    "4.Sterilization by filtration: microbial aerosols; air samplers; filter sizing; media sterilization.",
    "5.Kinetics and stoichiometry of microbial growth and product formation: transformation rates and conversion factors; classification of fermentative processes based on cell growth and product formation rates; influence of substrate concentration on cell growth rate; stoichiometry of microbial growth and product formation."
)
$r1.Text = ($part1 -join $lb)

# ---------------------------------------------------------------------------
# Block 2: "Critério" section text (grading formula) -> split into several
# lines separated by line breaks.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute(
    "A nota final (NF) será composta pelas médias M1  e M2,calculadas conforme segue:M1=P1+a1×T1M2=P2+a2×T2Em que:-P1 e P2 são as notas da primeira e da segunda prova escrita, respectivamente (notas de zero a dez).-T1 e T2 são as notas médias dos trabalhos (notas de zero a dez) realizados antes da primeira e da segunda prova escrita, respectivamente.-a1 e a2 são os fatores multiplicadores das notas médias dos trabalhos, a serem definidos pelo docente antes do início de cada turma com base nas atividades específicas a serem propostas. Os valores serão ≥0,1, sendo informados aos alunos no início do semestre. Em todos os casos, os valores máximos para M1 e M2 serão “dez”, sendo desconsideradas pontuações superiores.O cálculo de NF será feito conforme segue:NF=(M1+2×M2)/3Serão aprovados os alunos que obtiverem NF maior ou igual 5,0.",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0
) | Out-Null
$part2 = @(
    "A nota final (NF) será composta pelas médias M1  e M2,calculadas conforme segue:",
    "M1=P1+a1×T1",
    "M2=P2+a2×T2",
    "Em que:",
    "-P1 e P2 são as notas da primeira e da segunda prova escrita, respectivamente (notas de zero a dez).",
    "-T1 e T2 são as notas médias dos trabalhos (notas de zero a dez) realizados antes da primeira e da segunda prova escrita, respectivamente.",
    "-a1 e a2 são os fatores multiplicadores das notas médias dos trabalhos, a serem definidos pelo docente antes do início de cada turma com base nas atividades específicas a serem propostas. Os valores serão ≥0,1, sendo informados aos alunos no início do semestre. ",
    "Em todos os casos, os valores máximos para M1 e M2 serão “dez”, sendo desconsideradas pontuações superiores.",
    "O cálculo de NF será feito conforme segue:",
    "NF=(M1+2×M2)/3",
    "Serão aprovados os alunos que obtiverem NF maior ou igual 5,0."
)
$r2.Text = ($part2 -join $lb)

# ---------------------------------------------------------------------------
# Block 3: "Norma de recuperação" section text -> split into 3 lines
# separated by line breaks.
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute(
    "Será oferecido um programa de recuperação, sendo este avaliado por uma prova escrita final (PR). A média de recuperação (MR) será calculada conforme segue: MR=(NF+PR)/2Serão aprovados os alunos que obtiverem MR maior ou igual a 5,0.",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0
) | Out-Null
$part3 = @(
    "Será oferecido um programa de recuperação, sendo este avaliado por uma prova escrita final (PR). A média de recuperação (MR) será calculada conforme segue: ",
    "MR=(NF+PR)/2",
    "Serão aprovados os alunos que obtiverem MR maior ou igual a 5,0."
)
$r3.Text = ($part3 -join $lb)

# ---------------------------------------------------------------------------
# Block 4: Bibliography paragraph -> split into 5 separate references
# separated by line breaks.
# ---------------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute(
    "ALTERTHUM, F.; SCHMIDELL, W.; LIMA, U. A.; MORAES. M. O. Biotecnologia Industrial. Volume 1: Fundamentos. 2ª Edição. São Paulo: Blucher, 2020. ISBN 978-85-212-1897-5 (e-Book); 978-85-212-1898-2 (Impresso).ALTERTHUM, F.; SCHMIDELL, W.; LIMA, U. A.; MORAES. M. O. (Org.). Biotecnologia Industrial. Volume 2: Engenharia Bioquímica. 2ª Edição. São Paulo: Blucher, 2021. p. 37-52.  ISBN 978-65-5506-019-5 (e-Book); 978-65-5506-018-8 (Impresso).BORZANI, W. Processo Biotecnológico Industrial Genérico. In: BORZANI, W.; SCHMIDELL, W.; LIMA, U. A.; AQUARONE, E. Biotecnologia Industrial. Volume 1: Fundamentos. São Paulo: Editora Edgard Blücher Ltda, 2001. ISBN 978-85-212-0278-3.DORAN P.M.; MORRISSEY, K.; CARLSON, R. P. Bioprocess Engineering Principles, 3rd edition, Academic Press, 2024. ISBN 978-0128221914SHULER, M. L.; KARGI, F.; DELISA, M. Bioprocess Engineering: Basic Concepts (3rd Edition) (Prentice Hall International Series in the Physical and Chemical Engineering Sciences) 3rd Edition. Prentice Hall; 3 edition, 2017. ISBN: 978-0137062706",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0
) | Out-Null
$part4 = @(
    "ALTERTHUM, F.; SCHMIDELL, W.; LIMA, U. A.; MORAES. M. O. Biotecnologia Industrial. Volume 1: Fundamentos. 2ª Edição. São Paulo: Blucher, 2020. ISBN 978-85-212-1897-5 (e-Book); 978-85-212-1898-2 (Impresso).",
    "ALTERTHUM, F.; SCHMIDELL, W.; LIMA, U. A.; MORAES. M. O. (Org.). Biotecnologia Industrial. Volume 2: Engenharia Bioquímica. 2ª Edição. São Paulo: Blucher, 2021. p. 37-52.  ISBN 978-65-5506-019-5 (e-Book); 978-65-5506-018-8 (Impresso).",
    "BORZANI, W. Processo Biotecnológico Industrial Genérico. In: BORZANI, W.; SCHMIDELL, W.; LIMA, U. A.; AQUARONE, E. Biotecnologia Industrial. Volume 1: Fundamentos. São Paulo: Editora Edgard Blücher Ltda, 2001. ISBN 978-85-212-0278-3.",
    "DORAN P.M.; MORRISSEY, K.; CARLSON, R. P. Bioprocess Engineering Principles, 3rd edition, Academic Press, 2024. ISBN 978-0128221914",
    "SHULER, M. L.; KARGI, F.; DELISA, M. Bioprocess Engineering: Basic Concepts (3rd Edition) (Prentice Hall International Series in the Physical and Chemical Engineering Sciences) 3rd Edition. Prentice Hall; 3 edition, 2017. ISBN: 978-0137062706"
)
$r4.Text = ($part4 -join $lb)
